# Update the vm_pu.xlsx result values for the "case with 380 kV done" run.
# Columns C,D,E,F,J,K,L,M,N on rows 2-25 of Sheet1 hold per-bus voltage
# magnitudes (p.u.); column G is the slack bus and stays at 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.022350787414186
$ws.Range("D2").Value = 1.031528190178009
$ws.Range("E2").Value = 1.023070562848129
$ws.Range("F2").Value = 1.038231974351167
$ws.Range("J2").Value = 1.027537115733058
$ws.Range("K2").Value = 1.034336047918498
$ws.Range("L2").Value = 1.025903081956963
$ws.Range("M2").Value = 1.041020594774783
$ws.Range("N2").Value = 1.013110380410271

$ws.Range("C3").Value = 1.023593510720765
$ws.Range("D3").Value = 1.0327422910951
$ws.Range("E3").Value = 1.024135802793349
$ws.Range("F3").Value = 1.039598811691194
$ws.Range("J3").Value = 1.028416327729058
$ws.Range("K3").Value = 1.03535764479218
$ws.Range("L3").Value = 1.026774417607138
$ws.Range("M3").Value = 1.042195931920883
$ws.Range("N3").Value = 1.013413320091271

$ws.Range("C4").Value = 1.024397749068383
$ws.Range("D4").Value = 1.033528378173576
$ws.Range("E4").Value = 1.024825458674122
$ws.Range("F4").Value = 1.040483999279207
$ws.Range("J4").Value = 1.028984897215536
$ws.Range("K4").Value = 1.036018609460988
$ws.Range("L4").Value = 1.027338015034045
$ws.Range("M4").Value = 1.042956677041177
$ws.Range("N4").Value = 1.013608922337439

$ws.Range("C5").Value = 1.024735880523187
$ws.Range("D5").Value = 1.033858966674272
$ws.Range("E5").Value = 1.025115480973066
$ws.Range("F5").Value = 1.040856314753914
$ws.Range("J5").Value = 1.029223844021568
$ws.Range("K5").Value = 1.036296462350704
$ws.Range("L5").Value = 1.027574901192634
$ws.Range("M5").Value = 1.043276549872753
$ws.Range("N5").Value = 1.013691053268918

$ws.Range("C6").Value = 1.024792656095009
$ws.Range("D6").Value = 1.033914480885
$ws.Range("E6").Value = 1.025164182366333
$ws.Range("F6").Value = 1.040918838941032
$ws.Range("J6").Value = 1.029263959604413
$ws.Range("K6").Value = 1.036343114147181
$ws.Range("L6").Value = 1.027614672514299
$ws.Range("M6").Value = 1.04333026124927
$ws.Range("N6").Value = 1.013704837543325

$ws.Range("C7").Value = 1.024402267077738
$ws.Range("D7").Value = 1.033532795049188
$ws.Range("E7").Value = 1.024829333607175
$ws.Range("F7").Value = 1.040488973456156
$ws.Range("J7").Value = 1.028988090346214
$ws.Range("K7").Value = 1.036022322208186
$ws.Range("L7").Value = 1.027341180514156
$ws.Range("M7").Value = 1.042960950976578
$ws.Range("N7").Value = 1.013610020168529

$ws.Range("C8").Value = 1.022770749646192
$ws.Range("D8").Value = 1.031938402518974
$ws.Range("E8").Value = 1.023430488954062
$ws.Range("F8").Value = 1.038693749124722
$ws.Range("J8").Value = 1.027834320712148
$ws.Range("K8").Value = 1.034681318314149
$ws.Range("L8").Value = 1.026197599202938
$ws.Range("M8").Value = 1.041417760474974
$ws.Range("N8").Value = 1.01321284743053

$ws.Range("C9").Value = 1.019896561574956
$ws.Range("D9").Value = 1.029132487648931
$ws.Range("E9").Value = 1.020968353559982
$ws.Range("F9").Value = 1.035535986312455
$ws.Range("J9").Value = 1.025798565733092
$ws.Range("K9").Value = 1.032317632095826
$ws.Range("L9").Value = 1.024180759041069
$ws.Range("M9").Value = 1.038700078505805
$ws.Range("N9").Value = 1.012509748483676

$ws.Range("C10").Value = 1.017980783264892
$ws.Range("D10").Value = 1.027264179551022
$ws.Range("E10").Value = 1.019328731903349
$ws.Range("F10").Value = 1.033434454454585
$ws.Range("J10").Value = 1.0244395104461
$ws.Range("K10").Value = 1.030741283319735
$ws.Range("L10").Value = 1.022834971900181
$ws.Range("M10").Value = 1.036889244533728
$ws.Range("N10").Value = 1.012038824820527

$ws.Range("C11").Value = 1.017151276390657
$ws.Range("D11").Value = 1.026455699255192
$ws.Range("E11").Value = 1.018619165307329
$ws.Range("F11").Value = 1.032525296757998
$ws.Range("J11").Value = 1.023850556107359
$ws.Range("K11").Value = 1.030058549531312
$ws.Range("L11").Value = 1.022251921506253
$ws.Range("M11").Value = 1.036105329595163
$ws.Range("N11").Value = 1.011834385359114

$ws.Range("C12").Value = 1.016843163104841
$ws.Range("D12").Value = 1.026155467707971
$ws.Range("E12").Value = 1.018355659059484
$ws.Range("F12").Value = 1.032187715109719
$ws.Range("J12").Value = 1.02363171906525
$ws.Range("K12").Value = 1.029804924881122
$ws.Range("L12").Value = 1.022035301521377
$ws.Range("M12").Value = 1.035814174280222
$ws.Range("N12").Value = 1.011758367914563

$ws.Range("C13").Value = 1.016909254388389
$ws.Range("D13").Value = 1.02621986508416
$ws.Range("E13").Value = 1.018412179430913
$ws.Range("F13").Value = 1.032260122139005
$ws.Range("J13").Value = 1.023678663695606
$ws.Range("K13").Value = 1.029859329483643
$ws.Range("L13").Value = 1.022081769492735
$ws.Range("M13").Value = 1.035876626968145
$ws.Range("N13").Value = 1.011774677520252

$ws.Range("C14").Value = 1.017125807617683
$ws.Range("D14").Value = 1.026430880526107
$ws.Range("E14").Value = 1.018597382618046
$ws.Range("F14").Value = 1.032497389723697
$ws.Range("J14").Value = 1.023832468472985
$ws.Range("K14").Value = 1.030037585400715
$ws.Range("L14").Value = 1.022234016642251
$ws.Range("M14").Value = 1.036081262099261
$ws.Range("N14").Value = 1.011828103358284

$ws.Range("C15").Value = 1.017259233445825
$ws.Range("D15").Value = 1.026560903839274
$ws.Range("E15").Value = 1.018711500096753
$ws.Range("F15").Value = 1.032643593921437
$ws.Range("J15").Value = 1.023927222918935
$ws.Range("K15").Value = 1.030147411120228
$ws.Range("L15").Value = 1.022327814578219
$ws.Range("M15").Value = 1.036207347869417
$ws.Range("N15").Value = 1.011861010229249

$ws.Range("C16").Value = 1.018035835429233
$ws.Range("D16").Value = 1.027317846245128
$ws.Range("E16").Value = 1.019375831785874
$ws.Range("F16").Value = 1.033494809256067
$ws.Range("J16").Value = 1.024478587248081
$ws.Range("K16").Value = 1.030786590459599
$ws.Range("L16").Value = 1.022873660240158
$ws.Range("M16").Value = 1.036941274122437
$ws.Range("N16").Value = 1.012052381666033

$ws.Range("C17").Value = 1.018522985208828
$ws.Range("D17").Value = 1.027792790841299
$ws.Range("E17").Value = 1.019792655448253
$ws.Range("F17").Value = 1.034028971342203
$ws.Range("J17").Value = 1.024824315165395
$ws.Range("K17").Value = 1.031187485193935
$ws.Range("L17").Value = 1.023215969099819
$ws.Range("M17").Value = 1.037401695154321
$ws.Range("N17").Value = 1.012172282668454

$ws.Range("C18").Value = 1.018807135293721
$ws.Range("D18").Value = 1.028069867312775
$ws.Range("E18").Value = 1.020035820392533
$ws.Range("F18").Value = 1.034340618327948
$ws.Range("J18").Value = 1.025025926763864
$ws.Range("K18").Value = 1.031421304602492
$ws.Range("L18").Value = 1.023415601722137
$ws.Range("M18").Value = 1.037670269284661
$ws.Range("N18").Value = 1.012242168134054

$ws.Range("C19").Value = 1.018904023879935
$ws.Range("D19").Value = 1.028164351670755
$ws.Range("E19").Value = 1.020118739981219
$ws.Range("F19").Value = 1.0344468954371
$ws.Range("J19").Value = 1.025094663423481
$ws.Range("K19").Value = 1.031501028381117
$ws.Range("L19").Value = 1.023483666190937
$ws.Range("M19").Value = 1.037761849312173
$ws.Range("N19").Value = 1.012265988675049

$ws.Range("C20").Value = 1.018470718272737
$ws.Range("D20").Value = 1.027741828698423
$ws.Range("E20").Value = 1.01974793021937
$ws.Range("F20").Value = 1.033971652629332
$ws.Range("J20").Value = 1.024787226567877
$ws.Range("K20").Value = 1.031144474624801
$ws.Range("L20").Value = 1.023179245750419
$ws.Range("M20").Value = 1.037352294459242
$ws.Range("N20").Value = 1.012159423679525

$ws.Range("C21").Value = 1.01706203804557
$ws.Range("D21").Value = 1.02636873971436
$ws.Range("E21").Value = 1.018542843298231
$ws.Range("F21").Value = 1.032427517046178
$ws.Range("J21").Value = 1.023787178834855
$ws.Range("K21").Value = 1.029985094237149
$ws.Range("L21").Value = 1.022189185020519
$ws.Range("M21").Value = 1.036021001468365
$ws.Range("N21").Value = 1.011812372979428

$ws.Range("C22").Value = 1.016176357773045
$ws.Range("D22").Value = 1.025505851257195
$ws.Range("E22").Value = 1.017785494204068
$ws.Range("F22").Value = 1.031457350459447
$ws.Range("J22").Value = 1.023157984808569
$ws.Range("K22").Value = 1.029255988901635
$ws.Range("L22").Value = 1.021566409584932
$ws.Range("M22").Value = 1.035184111500782
$ws.Range("N22").Value = 1.011593708098701

$ws.Range("C23").Value = 1.016645873480442
$ws.Range("D23").Value = 1.025963245108554
$ws.Range("E23").Value = 1.018186947949124
$ws.Range("F23").Value = 1.03197158949945
$ws.Range("J23").Value = 1.023491573321971
$ws.Range("K23").Value = 1.029642517048552
$ws.Range("L23").Value = 1.021896582198739
$ws.Range("M23").Value = 1.035627749684604
$ws.Range("N23").Value = 1.011709670270188

$ws.Range("C24").Value = 1.018494335439357
$ws.Range("D24").Value = 1.027764856144744
$ws.Range("E24").Value = 1.019768139505433
$ws.Range("F24").Value = 1.033997552243424
$ws.Range("J24").Value = 1.024803985449674
$ws.Range("K24").Value = 1.031163909298148
$ws.Range("L24").Value = 1.023195839545054
$ws.Range("M24").Value = 1.037374616450248
$ws.Range("N24").Value = 1.012165234259987

$ws.Range("C25").Value = 1.020639535189353
$ws.Range("D25").Value = 1.02985746907419
$ws.Range("E25").Value = 1.02160455068109
$ws.Range("F25").Value = 1.036351690044639
$ws.Range("J25").Value = 1.032928791908624
$ws.Range("K25").Value = 1.039402486532032
$ws.Range("L25").Value = 1.024702370860777
$ws.Range("M25").Value = 1.01269190116642
